$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the header row (subject numbers) for columns B:E
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Update CON row (row 2) values for columns B:E
$ws.Range("B2").Value = 35.00680295715339
$ws.Range("C2").Value = 45.407415285381461
$ws.Range("D2").Value = 38.816475182433159
$ws.Range("E2").Value = 43.068835150836037

# Update STR row (row 3) values for columns B:E
$ws.Range("B3").Value = 43.216688876332171
$ws.Range("C3").Value = 48.834862836497692
$ws.Range("D3").Value = 45.18835969066658
$ws.Range("E3").Value = 42.78755270190009

# Shrink the active selection to match the new working range
$ws.Range("B1:E3").Select()
